# Update the "g18.7" data sheet: extend the time series for
# Brasil, Nordeste and Sergipe from 2010-2021 to 2010-2022.
# This shifts every row from A15:C37 downward (inserting the new
# "2022" observation for each region after its prior last year) and
# appends a brand-new "Sergipe / 2022" row at the bottom, growing the
# sheet from A1:C37 to A1:C40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset (header stays in row 1; data starts row 2).
$data = @(
    ,@("Brasil", "01/01/2010", 0.3154941217960625)
    ,@("Brasil", "01/01/2011", 0.3155312903960848)
    ,@("Brasil", "01/01/2012", 0.3108866536595493)
    ,@("Brasil", "01/01/2013", 0.310993830647957)
    ,@("Brasil", "01/01/2014", 0.3056745780905651)
    ,@("Brasil", "01/01/2015", 0.293626307380059)
    ,@("Brasil", "01/01/2016", 0.2892400683052905)
    ,@("Brasil", "01/01/2017", 0.2862688760429477)
    ,@("Brasil", "01/01/2018", 0.2822170017670938)
    ,@("Brasil", "01/01/2019", 0.2778353668508045)
    ,@("Brasil", "01/01/2020", 0.2445637434267342)
    ,@("Brasil", "01/01/2021", 0.250807965922708)
    ,@("Brasil", "01/01/2022", 0.2626141501055874)
    ,@("Nordeste", "01/01/2010", 0.3447120207112737)
    ,@("Nordeste", "01/01/2011", 0.3447559014477932)
    ,@("Nordeste", "01/01/2012", 0.3337083375505088)
    ,@("Nordeste", "01/01/2013", 0.343725831195954)
    ,@("Nordeste", "01/01/2014", 0.3333886421255093)
    ,@("Nordeste", "01/01/2015", 0.3233650204064076)
    ,@("Nordeste", "01/01/2016", 0.3232012645743133)
    ,@("Nordeste", "01/01/2017", 0.3233904122747474)
    ,@("Nordeste", "01/01/2018", 0.3182142809705919)
    ,@("Nordeste", "01/01/2019", 0.3089531582700869)
    ,@("Nordeste", "01/01/2020", 0.2781971063295547)
    ,@("Nordeste", "01/01/2021", 0.2805786372078377)
    ,@("Nordeste", "01/01/2022", 0.2907529924918893)
    ,@("Sergipe", "01/01/2010", 0.3367492251417083)
    ,@("Sergipe", "01/01/2011", 0.3697077440467428)
    ,@("Sergipe", "01/01/2012", 0.3655982881532424)
    ,@("Sergipe", "01/01/2013", 0.339706937660594)
    ,@("Sergipe", "01/01/2014", 0.3659540193657965)
    ,@("Sergipe", "01/01/2015", 0.3600606909430439)
    ,@("Sergipe", "01/01/2016", 0.3334233419819904)
    ,@("Sergipe", "01/01/2017", 0.3442066420664207)
    ,@("Sergipe", "01/01/2018", 0.3806352549432577)
    ,@("Sergipe", "01/01/2019", 0.3691164715286594)
    ,@("Sergipe", "01/01/2020", 0.3348227443675505)
    ,@("Sergipe", "01/01/2021", 0.2995685976776077)
    ,@("Sergipe", "01/01/2022", 0.3415898922316227)
)

$lastRow = 1 + $data.Count

# Make sure column B keeps storing the dates as plain text (as in the
# original file) instead of letting Excel auto-convert "dd/mm/yyyy"
# strings into date serial numbers.
$ws.Range("B2:B$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
}
